# Natmi following Dr Hou advice
#
# The ligand-receptor signalling table (Rarres2-Ccrl2) is recomputed so that
# the "Target cluster" column now also includes the FAPs cluster (previously
# missing), turning the Sending x Target matrix into a full 3x3 grid over
# {ECs, FAPs, sCs}, and all of the expression/specificity statistics are
# refreshed to match the new cell counts (3 cells per cluster instead of 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rarres2"
$ws.Range("C2").Value = "Ccrl2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.916729333333334
$ws.Range("H2").Value = 11.750188
$ws.Range("I2").Value = 0.05842616646845182
$ws.Range("J2").Value = 0.05842616646845181
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.007593666666667
$ws.Range("N2").Value = 6.022781
$ws.Range("O2").Value = 0.03411486246698893
$ws.Range("P2").Value = 0.03411486246698894
$ws.Range("Q2").Value = 7.863201003647557
$ws.Range("R2").Value = 70.76880903282802
$ws.Range("S2").Value = 0.001993200633544634
$ws.Range("T2").Value = 0.001993200633544634

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rarres2"
$ws.Range("C3").Value = "Ccrl2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.916729333333334
$ws.Range("H3").Value = 11.750188
$ws.Range("I3").Value = 0.05842616646845182
$ws.Range("J3").Value = 0.05842616646845181
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.464105
$ws.Range("N3").Value = 1.392315
$ws.Range("O3").Value = 0.007886495413950083
$ws.Range("P3").Value = 0.007886495413950083
$ws.Range("Q3").Value = 1.817773667246667
$ws.Range("R3").Value = 16.35996300522
$ws.Range("S3").Value = 0.0004607776939081294
$ws.Range("T3").Value = 0.0004607776939081294

# Row 4: ECs -> M2
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Rarres2"
$ws.Range("C4").Value = "Ccrl2"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.916729333333334
$ws.Range("H4").Value = 11.750188
$ws.Range("I4").Value = 0.05842616646845182
$ws.Range("J4").Value = 0.05842616646845181
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 54.02549366666667
$ws.Range("N4").Value = 162.076481
$ws.Range("O4").Value = 0.9180504584922721
$ws.Range("P4").Value = 0.9180504584922721
$ws.Range("Q4").Value = 211.6032357920476
$ws.Range("R4").Value = 1904.429122128428
$ws.Range("S4").Value = 0.053638168914308
$ws.Range("T4").Value = 0.053638168914308

# Row 5: ECs -> sCs
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Rarres2"
$ws.Range("C5").Value = "Ccrl2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.916729333333334
$ws.Range("H5").Value = 11.750188
$ws.Range("I5").Value = 0.05842616646845182
$ws.Range("J5").Value = 0.05842616646845181
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.350873333333333
$ws.Range("N5").Value = 7.052619999999999
$ws.Range("O5").Value = 0.03994818362678892
$ws.Range("P5").Value = 0.03994818362678893
$ws.Range("Q5").Value = 9.207734543617777
$ws.Range("R5").Value = 82.86961089256
$ws.Range("S5").Value = 0.002334019226691051
$ws.Range("T5").Value = 0.002334019226691051

# Row 6: FAPs -> ECs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Rarres2"
$ws.Range("C6").Value = "Ccrl2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 48.53546666666667
$ws.Range("H6").Value = 145.6064
$ws.Range("I6").Value = 0.7240074597335789
$ws.Range("J6").Value = 0.7240074597335789
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.007593666666667
$ws.Range("N6").Value = 6.022781
$ws.Range("O6").Value = 0.03411486246698893
$ws.Range("P6").Value = 0.03411486246698894
$ws.Range("Q6").Value = 97.43949548871112
$ws.Range("R6").Value = 876.9554593984001
$ws.Range("S6").Value = 0.02469941491388507
$ws.Range("T6").Value = 0.02469941491388507

# Row 7: FAPs -> FAPs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Rarres2"
$ws.Range("C7").Value = "Ccrl2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 48.53546666666667
$ws.Range("H7").Value = 145.6064
$ws.Range("I7").Value = 0.7240074597335789
$ws.Range("J7").Value = 0.7240074597335789
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.464105
$ws.Range("N7").Value = 1.392315
$ws.Range("O7").Value = 0.007886495413950083
$ws.Range("P7").Value = 0.007886495413950083
$ws.Range("Q7").Value = 22.52555275733334
$ws.Range("R7").Value = 202.729974816
$ws.Range("S7").Value = 0.005709881510854519
$ws.Range("T7").Value = 0.005709881510854519

# Row 8: FAPs -> M2
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Rarres2"
$ws.Range("C8").Value = "Ccrl2"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 48.53546666666667
$ws.Range("H8").Value = 145.6064
$ws.Range("I8").Value = 0.7240074597335789
$ws.Range("J8").Value = 0.7240074597335789
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 54.02549366666667
$ws.Range("N8").Value = 162.076481
$ws.Range("O8").Value = 0.9180504584922721
$ws.Range("P8").Value = 0.9180504584922721
$ws.Range("Q8").Value = 2622.152547008712
$ws.Range("R8").Value = 23599.3729230784
$ws.Range("S8").Value = 0.6646753803602373
$ws.Range("T8").Value = 0.6646753803602373

# Row 9: FAPs -> sCs
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Rarres2"
$ws.Range("C9").Value = "Ccrl2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 48.53546666666667
$ws.Range("H9").Value = 145.6064
$ws.Range("I9").Value = 0.7240074597335789
$ws.Range("J9").Value = 0.7240074597335789
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.350873333333333
$ws.Range("N9").Value = 7.052619999999999
$ws.Range("O9").Value = 0.03994818362678892
$ws.Range("P9").Value = 0.03994818362678893
$ws.Range("Q9").Value = 114.1007343075555
$ws.Range("R9").Value = 1026.906608768
$ws.Range("S9").Value = 0.028922782948602
$ws.Range("T9").Value = 0.028922782948602

# Row 10: sCs -> ECs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Rarres2"
$ws.Range("C10").Value = "Ccrl2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 14.58505066666667
$ws.Range("H10").Value = 43.755152
$ws.Range("I10").Value = 0.2175663737979692
$ws.Range("J10").Value = 0.2175663737979692
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.007593666666667
$ws.Range("N10").Value = 6.022781
$ws.Range("O10").Value = 0.03411486246698893
$ws.Range("P10").Value = 0.03411486246698894
$ws.Range("Q10").Value = 29.28085534641245
$ws.Range("R10").Value = 263.527698117712
$ws.Range("S10").Value = 0.007422246919559224
$ws.Range("T10").Value = 0.007422246919559225

# Row 11: sCs -> FAPs
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Rarres2"
$ws.Range("C11").Value = "Ccrl2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 14.58505066666667
$ws.Range("H11").Value = 43.755152
$ws.Range("I11").Value = 0.2175663737979692
$ws.Range("J11").Value = 0.2175663737979692
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.464105
$ws.Range("N11").Value = 1.392315
$ws.Range("O11").Value = 0.007886495413950083
$ws.Range("P11").Value = 0.007886495413950083
$ws.Range("Q11").Value = 6.768994939653334
$ws.Range("R11").Value = 60.92095445688
$ws.Range("S11").Value = 0.001715836209187434
$ws.Range("T11").Value = 0.001715836209187434

# Row 12: sCs -> M2
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Rarres2"
$ws.Range("C12").Value = "Ccrl2"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 14.58505066666667
$ws.Range("H12").Value = 43.755152
$ws.Range("I12").Value = 0.2175663737979692
$ws.Range("J12").Value = 0.2175663737979692
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 54.02549366666667
$ws.Range("N12").Value = 162.076481
$ws.Range("O12").Value = 0.9180504584922721
$ws.Range("P12").Value = 0.9180504584922721
$ws.Range("Q12").Value = 787.9645624200125
$ws.Range("R12").Value = 7091.681061780112
$ws.Range("S12").Value = 0.1997369092177267
$ws.Range("T12").Value = 0.1997369092177267

# Row 13: sCs -> sCs
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Rarres2"
$ws.Range("C13").Value = "Ccrl2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 14.58505066666667
$ws.Range("H13").Value = 43.755152
$ws.Range("I13").Value = 0.2175663737979692
$ws.Range("J13").Value = 0.2175663737979692
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.350873333333333
$ws.Range("N13").Value = 7.052619999999999
$ws.Range("O13").Value = 0.03994818362678892
$ws.Range("P13").Value = 0.03994818362678893
$ws.Range("Q13").Value = 34.28760667758222
$ws.Range("R13").Value = 308.58846009824
$ws.Range("S13").Value = 0.008691381451495873
$ws.Range("T13").Value = 0.008691381451495875
